# Update the daily report: one extra day is inserted on 2021-02-08
# (date serial 44235), pushing every following row down by one, and a
# further day is appended at the very end on 2021-03-02 (serial
# 44257). Column C (7-day rolling sum) and column D (per-100k figure)
# are refreshed for every row whose rolling window changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 114 and 115 do not exist yet. Clone the formatting (date number
# format/border/font on col A, plain blank cells on B:D) from row 113
# so the new rows look exactly like their neighbours before we fill in
# their data.
$formatSource = $ws.Range("A113:D113")
$formatSource.Copy()
$ws.Range("A114:D115").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Final (row, date-serial, nuovi pos., somma mobile 7gg., somma mobile per 100mila)
# for every row that differs from the original sheet (rows 90-115).
$rows = @(
    @(90, 44232, 2, 8, 91.99632014719411),
    @(91, 44233, 2, 7, 80.49678012879485),
    @(92, 44234, 2, 7, 80.49678012879485),
    @(93, 44235, 1, 10, 114.9954001839926),
    @(94, 44236, 0, 15, 172.493100275989),
    @(95, 44237, 0, 13, 149.4940202391904),
    @(96, 44238, 3, 14, 160.9935602575897),
    @(97, 44239, 7, 19, 218.491260349586),
    @(98, 44240, 0, 19, 218.491260349586),
    @(99, 44241, 3, 19, 218.491260349586),
    @(100, 44242, 6, 19, 218.491260349586),
    @(101, 44243, 0, 15, 172.493100275989),
    @(102, 44244, 0, 15, 172.493100275989),
    @(103, 44245, 3, 16, 183.9926402943882),
    @(104, 44246, 3, 11, 126.4949402023919),
    @(105, 44247, 0, 16, 183.9926402943882),
    @(106, 44248, 4, 17, 195.4921803127875),
    @(107, 44249, 1, 21, 241.4903403863846),
    @(108, 44250, 5, 21, 241.4903403863846),
    @(109, 44251, 1, 25, 287.4885004599816),
    @(110, 44252, 7, 26, 298.9880404783809),
    @(111, 44253, 3, 34, 390.984360625575),
    @(112, 44254, 4, 29, 333.4866605335786)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value2 = $row[1]
    $ws.Cells.Item($r, 2).Value2 = $row[2]
    $ws.Cells.Item($r, 3).Value2 = $row[3]
    $ws.Cells.Item($r, 4).Value2 = $row[4]
}

# Rows 113-115 only carry a date and a "nuovi pos." count; columns C
# and D stay blank (no 7-day window available yet), matching the
# existing trailing rows.
$tailRows = @(
    @(113, 44255, 5),
    @(114, 44256, 9),
    @(115, 44257, 0)
)

foreach ($row in $tailRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value2 = $row[1]
    $ws.Cells.Item($r, 2).Value2 = $row[2]
}
